# Change the algorithm to replace cells inside tables.
#
# The new table-cell replacement algorithm no longer needs a second,
# separate "donor" table (title {{table2}}) to source extra blank cells
# from - it now replaces/grows cells directly inside the single target
# table. As a visible side effect on this fixture slide:
#   - the {{table1}} table is repositioned (it now sits lower on the
#     slide, since the layout no longer has to leave room for a second
#     table underneath it);
#   - the now-unused {{table2}} helper table is removed entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$table1Shape = $null
$table2Shape = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Title -eq "{{table1}}") {
        $table1Shape = $sh
    } elseif ($sh.Title -eq "{{table2}}") {
        $table2Shape = $sh
    }
}

# Move {{table1}} down on the slide (457200, 404664) EMU -> (457200, 1304776) EMU.
# Shape.Top/Left are expressed in points (1 pt = 12700 EMU).
if ($table1Shape -ne $null) {
    $table1Shape.Left = 36
    $table1Shape.Top = 102.73827
}

# Drop the now-obsolete {{table2}} helper table altogether.
if ($table2Shape -ne $null) {
    $table2Shape.Delete()
}
